# Add two new rows to the phi_constants sheet:
#  - "Resources [of Hydro]" right after the "Hydro" row (row 25 -> new row 26)
#  - "Resources [of Primary solid biofuels]" right after the "Primary solid
#    biofuels" row (originally row 46, now row 47 after the first insert -> new row 48)
#
# Both new rows are inserted the way Excel's UI "Insert Copied Cells" works:
# copy the row immediately above, insert it (shifting everything below down by
# one row, carrying over formatting), then adjust the specific cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) "Resources [of Hydro]" after "Hydro" (row 25) ---------------------
$ws.Rows.Item(25).Copy()
$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value = "Resources [of Hydro]"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = $false

# --- 2) "Resources [of Primary solid biofuels]" after "Primary solid
#        biofuels" (originally row 46, now row 47 after the first insert) --
$ws.Rows.Item(47).Copy()
$ws.Rows.Item(48).Insert()

$ws.Range("A48").Value = "Resources [of Primary solid biofuels]"
$ws.Range("B48").Formula = "=Serrenho_2013!B7"
$ws.Range("C48").Value = $false
$ws.Range("D48").Value = "Serrenho et al (2013)"

# --- Restore the selection to mirror the committed workbook's last state ---
$ws.Range("B47:D48").Select()

$wb.Save()
